$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Room Range values in column A to include wing prefix
$ws.Range("A6").Value = "B401-B416"
$ws.Range("A7").Value = "B422-B433"
$ws.Range("A8").Value = "B401-B416"
$ws.Range("A9").Value = "B422-B433"
$ws.Range("A4").Value = "A201-A234"

# Clear the now-redundant Wing column (D2:D9)
$ws.Range("D2:D9").ClearContents()

$ws.Range("D2").Select()
